$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record (row 132) was inserted into the daily log, pushing the
# existing rows 132:198 down to 133:199.
$ws.Rows(132).Insert()

# Populate the newly inserted row with its data.
$ws.Cells.Item(132, 1).Value = 5
$ws.Cells.Item(132, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(132, 3).Value = "Maule"
$ws.Cells.Item(132, 4).Value = 44529
$ws.Cells.Item(132, 5).Value = 7
$ws.Cells.Item(132, 6).Value = "Fruta"
$ws.Cells.Item(132, 7).Value = 100103
$ws.Cells.Item(132, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(132, 9).Value = 100103004
$ws.Cells.Item(132, 10).Value = "Durazno"
$ws.Cells.Item(132, 11).Value = "Florida King"
$ws.Cells.Item(132, 12).Value = "Especial"
$ws.Cells.Item(132, 13).Value = 200
$ws.Cells.Item(132, 14).Value = 10000
$ws.Cells.Item(132, 15).Value = 10000
$ws.Cells.Item(132, 16).Value = 10000
$ws.Cells.Item(132, 17).Value = "`$/caja 12 kilos empedrada"
$ws.Cells.Item(132, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(132, 19).Value = 833
$ws.Cells.Item(132, 20).Value = 12
